$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Prefix with an apostrophe so Excel always stores this as literal text
    # (shared string), never auto-coercing to Boolean/Date/Number.
    $range.Value = "'" + $text
}

$commit = "58784ee182f52ffb6c3c60040f6983b7fa51718e"
$newFile = "f607f7c6-f204-4dde-82d2-c1cadd15fb77.md"
$oldFile = "f77d4477-3f5e-44b0-9f38-0d86a46c8db3.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A..G
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item(2).Insert()
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

Set-TextValue $ws.Range("A2") $newFile
Set-TextValue $ws.Range("B2") ("e2e\" + $newFile)
Set-TextValue $ws.Range("C2") ".md"
Set-TextValue $ws.Range("D2") ""
Set-TextValue $ws.Range("E2") "Ready for handoff"
Set-TextValue $ws.Range("F2") "Ready for handoff"
Set-TextValue $ws.Range("G2") "2016-08-23 14:47:30"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Set-TextValue $ws.Range("A3") $oldFile
Set-TextValue $ws.Range("B3") ("e2e\" + $oldFile)
Set-TextValue $ws.Range("C3") ".md"
Set-TextValue $ws.Range("D3") ""
Set-TextValue $ws.Range("E3") "Ready for handoff"
Set-TextValue $ws.Range("F3") "Ready for handoff"
Set-TextValue $ws.Range("G3") "2016-08-23 14:46:56"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile", "", "", "e2e\$newFile") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$oldFile", "", "", "e2e\$oldFile") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A..P
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

Set-TextValue $ws.Range("A2") $newFile
Set-TextValue $ws.Range("B2") ".md"
Set-TextValue $ws.Range("C2") "Ready for handoff"
Set-TextValue $ws.Range("D2") "e2e"
Set-TextValue $ws.Range("E2") "ht"
Set-TextValue $ws.Range("F2") "False"
Set-TextValue $ws.Range("G2") "f607f7c6-f204-4dde-82d2-c1cadd15fb77.6fe5a1463f5ee4c99bd7d3af2e2db969f44e613b.zh-cn.xlf"
Set-TextValue $ws.Range("H2") "2016-08-23 14:47:24"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("I2") ""
Set-TextValue $ws.Range("J2") ""
Set-TextValue $ws.Range("K2") "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("L2") ""
Set-TextValue $ws.Range("M2") "True"
Set-TextValue $ws.Range("N2") ""
Set-TextValue $ws.Range("O2") "False"
Set-TextValue $ws.Range("P2") ""

Set-TextValue $ws.Range("A3") $oldFile
Set-TextValue $ws.Range("B3") ".md"
Set-TextValue $ws.Range("C3") "Ready for handoff"
Set-TextValue $ws.Range("D3") "e2e"
Set-TextValue $ws.Range("E3") "ht"
Set-TextValue $ws.Range("F3") "False"
Set-TextValue $ws.Range("G3") "f77d4477-3f5e-44b0-9f38-0d86a46c8db3.a8a84a9e223c61a77591c106a1c01914d0c00c38.zh-cn.xlf"
Set-TextValue $ws.Range("H3") "2016-08-23 14:46:52"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("I3") ""
Set-TextValue $ws.Range("J3") ""
Set-TextValue $ws.Range("K3") "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("L3") ""
Set-TextValue $ws.Range("M3") "True"
Set-TextValue $ws.Range("N3") ""
Set-TextValue $ws.Range("O3") "False"
Set-TextValue $ws.Range("P3") ""

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile", "", "", $newFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$oldFile", "", "", $oldFile) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A..P
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

Set-TextValue $ws.Range("A2") $newFile
Set-TextValue $ws.Range("B2") ".md"
Set-TextValue $ws.Range("C2") "Ready for handoff"
Set-TextValue $ws.Range("D2") "e2e"
Set-TextValue $ws.Range("E2") "ht"
Set-TextValue $ws.Range("F2") "False"
Set-TextValue $ws.Range("G2") "f607f7c6-f204-4dde-82d2-c1cadd15fb77.6fe5a1463f5ee4c99bd7d3af2e2db969f44e613b.de-de.xlf"
Set-TextValue $ws.Range("H2") "2016-08-23 14:47:30"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("I2") ""
Set-TextValue $ws.Range("J2") ""
Set-TextValue $ws.Range("K2") "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("L2") ""
Set-TextValue $ws.Range("M2") "True"
Set-TextValue $ws.Range("N2") ""
Set-TextValue $ws.Range("O2") "False"
Set-TextValue $ws.Range("P2") ""

Set-TextValue $ws.Range("A3") $oldFile
Set-TextValue $ws.Range("B3") ".md"
Set-TextValue $ws.Range("C3") "Ready for handoff"
Set-TextValue $ws.Range("D3") "e2e"
Set-TextValue $ws.Range("E3") "ht"
Set-TextValue $ws.Range("F3") "False"
Set-TextValue $ws.Range("G3") "f77d4477-3f5e-44b0-9f38-0d86a46c8db3.a8a84a9e223c61a77591c106a1c01914d0c00c38.de-de.xlf"
Set-TextValue $ws.Range("H3") "2016-08-23 14:46:56"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("I3") ""
Set-TextValue $ws.Range("J3") ""
Set-TextValue $ws.Range("K3") "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-TextValue $ws.Range("L3") ""
Set-TextValue $ws.Range("M3") "True"
Set-TextValue $ws.Range("N3") ""
Set-TextValue $ws.Range("O3") "False"
Set-TextValue $ws.Range("P3") ""

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$newFile", "", "", $newFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$oldFile", "", "", $oldFile) | Out-Null
